$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 118.0346986666667
$ws.Range("H2").Value = 354.104096
$ws.Range("I2").Value = 0.2666057129183408
$ws.Range("J2").Value = 0.2666057129183408
$ws.Range("O2").Value = 0.2013489143128838
$ws.Range("P2").Value = 0.2013489143128838
$ws.Range("Q2").Value = 148.4554274499307
$ws.Range("R2").Value = 1336.098847049376
$ws.Range("S2").Value = 0.0536807708457203
$ws.Range("T2").Value = 0.0536807708457203

# Row 3
$ws.Range("G3").Value = 118.0346986666667
$ws.Range("H3").Value = 354.104096
$ws.Range("I3").Value = 0.2666057129183408
$ws.Range("J3").Value = 0.2666057129183408
$ws.Range("M3").Value = 0.246708
$ws.Range("N3").Value = 0.740124
$ws.Range("O3").Value = 0.03949536580856015
$ws.Range("P3").Value = 0.03949536580856015
$ws.Range("Q3").Value = 29.120104438656
$ws.Range("R3").Value = 262.080939947904
$ws.Range("S3").Value = 0.01052969015836184
$ws.Range("T3").Value = 0.01052969015836184

# Row 4
$ws.Range("G4").Value = 118.0346986666667
$ws.Range("H4").Value = 354.104096
$ws.Range("I4").Value = 0.2666057129183408
$ws.Range("J4").Value = 0.2666057129183408
$ws.Range("M4").Value = 4.74207
$ws.Range("N4").Value = 14.22621
$ws.Range("O4").Value = 0.7591557198785561
$ws.Range("P4").Value = 0.759155719878556
$ws.Range("Q4").Value = 559.72880350624
$ws.Range("R4").Value = 5037.55923155616
$ws.Range("S4").Value = 0.2023952519142587
$ws.Range("T4").Value = 0.2023952519142587

# Row 5
$ws.Range("I5").Value = 0.4881754016778185
$ws.Range("J5").Value = 0.4881754016778186
$ws.Range("O5").Value = 0.2013489143128838
$ws.Range("P5").Value = 0.2013489143128838
$ws.Range("S5").Value = 0.09829358712208469
$ws.Range("T5").Value = 0.09829358712208471

# Row 6
$ws.Range("I6").Value = 0.4881754016778185
$ws.Range("J6").Value = 0.4881754016778186
$ws.Range("M6").Value = 0.246708
$ws.Range("N6").Value = 0.740124
$ws.Range("O6").Value = 0.03949536580856015
$ws.Range("P6").Value = 0.03949536580856015
$ws.Range("Q6").Value = 53.321133015612
$ws.Range("R6").Value = 479.890197140508
$ws.Range("S6").Value = 0.01928066606800623
$ws.Range("T6").Value = 0.01928066606800623

# Row 7
$ws.Range("I7").Value = 0.4881754016778185
$ws.Range("J7").Value = 0.4881754016778186
$ws.Range("M7").Value = 4.74207
$ws.Range("N7").Value = 14.22621
$ws.Range("O7").Value = 0.7591557198785561
$ws.Range("P7").Value = 0.759155719878556
$ws.Range("Q7").Value = 1024.90614507573
$ws.Range("R7").Value = 9224.15530568157
$ws.Range("S7").Value = 0.3706011484877276
$ws.Range("T7").Value = 0.3706011484877276

# Row 8
$ws.Range("G8").Value = 45.876452
$ws.Range("H8").Value = 137.629356
$ws.Range("I8").Value = 0.1036214293744632
$ws.Range("J8").Value = 0.1036214293744632
$ws.Range("O8").Value = 0.2013489143128838
$ws.Range("P8").Value = 0.2013489143128838
$ws.Range("Q8").Value = 57.700052344604
$ws.Range("R8").Value = 519.3004711014361
$ws.Range("S8").Value = 0.02086406230409732
$ws.Range("T8").Value = 0.02086406230409733

# Row 9
$ws.Range("G9").Value = 45.876452
$ws.Range("H9").Value = 137.629356
$ws.Range("I9").Value = 0.1036214293744632
$ws.Range("J9").Value = 0.1036214293744632
$ws.Range("M9").Value = 0.246708
$ws.Range("N9").Value = 0.740124
$ws.Range("O9").Value = 0.03949536580856015
$ws.Range("P9").Value = 0.03949536580856015
$ws.Range("Q9").Value = 11.318087720016
$ws.Range("R9").Value = 101.862789480144
$ws.Range("S9").Value = 0.004092566258750303
$ws.Range("T9").Value = 0.004092566258750304

# Row 10
$ws.Range("G10").Value = 45.876452
$ws.Range("H10").Value = 137.629356
$ws.Range("I10").Value = 0.1036214293744632
$ws.Range("J10").Value = 0.1036214293744632
$ws.Range("M10").Value = 4.74207
$ws.Range("N10").Value = 14.22621
$ws.Range("O10").Value = 0.7591557198785561
$ws.Range("P10").Value = 0.759155719878556
$ws.Range("Q10").Value = 217.54934673564
$ws.Range("R10").Value = 1957.94412062076
$ws.Range("S10").Value = 0.07866480081161556
$ws.Range("T10").Value = 0.07866480081161556

# Row 11
$ws.Range("G11").Value = 62.68962833333333
$ws.Range("H11").Value = 188.068885
$ws.Range("I11").Value = 0.1415974560293775
$ws.Range("J11").Value = 0.1415974560293775
$ws.Range("O11").Value = 0.2013489143128838
$ws.Range("P11").Value = 0.2013489143128838
$ws.Range("Q11").Value = 78.84643817479834
$ws.Range("R11").Value = 709.617943573185
$ws.Range("S11").Value = 0.02851049404098145
$ws.Range("T11").Value = 0.02851049404098145

# Row 12
$ws.Range("G12").Value = 62.68962833333333
$ws.Range("H12").Value = 188.068885
$ws.Range("I12").Value = 0.1415974560293775
$ws.Range("J12").Value = 0.1415974560293775
$ws.Range("M12").Value = 0.246708
$ws.Range("N12").Value = 0.740124
$ws.Range("O12").Value = 0.03949536580856015
$ws.Range("P12").Value = 0.03949536580856015
$ws.Range("Q12").Value = 15.46603282686
$ws.Range("R12").Value = 139.19429544174
$ws.Range("S12").Value = 0.005592443323441773
$ws.Range("T12").Value = 0.005592443323441774

# Row 13
$ws.Range("G13").Value = 62.68962833333333
$ws.Range("H13").Value = 188.068885
$ws.Range("I13").Value = 0.1415974560293775
$ws.Range("J13").Value = 0.1415974560293775
$ws.Range("M13").Value = 4.74207
$ws.Range("N13").Value = 14.22621
$ws.Range("O13").Value = 0.7591557198785561
$ws.Range("P13").Value = 0.759155719878556
$ws.Range("Q13").Value = 297.27860583065
$ws.Range("R13").Value = 2675.50745247585
$ws.Range("S13").Value = 0.1074945186649542
$ws.Range("T13").Value = 0.1074945186649542
